# Sync attendance_reports: normalize "Recorded By" ordering in column G.
# Every cell in column G whose value is exactly "dnasr281@gmail.com, System"
# is rewritten to "System, dnasr281@gmail.com" (order swapped), leaving all
# other cells (formatting, other values, etc.) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
